$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.342908333333333
$ws.Range("H2").Value = 4.028725
$ws.Range("I2").Value = 0.2879023314891748
$ws.Range("J2").Value = 0.2879023314891748
$ws.Range("M2").Value = 8.333446333333333
$ws.Range("N2").Value = 25.000339
$ws.Range("O2").Value = 0.3294294409523786
$ws.Range("P2").Value = 0.3294294409523787
$ws.Range("Q2").Value = 11.19105452641944
$ws.Range("R2").Value = 100.719490737775
$ws.Range("S2").Value = 0.09484350411136525
$ws.Range("T2").Value = 0.09484350411136526
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.342908333333333
$ws.Range("H3").Value = 4.028725
$ws.Range("I3").Value = 0.2879023314891748
$ws.Range("J3").Value = 0.2879023314891748
$ws.Range("O3").Value = 0.357886883212021
$ws.Range("P3").Value = 0.357886883212021
$ws.Range("Q3").Value = 12.15778290105833
$ws.Range("R3").Value = 109.420046109525
$ws.Range("S3").Value = 0.1030364680861348
$ws.Range("T3").Value = 0.1030364680861348
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.342908333333333
$ws.Range("H4").Value = 4.028725
$ws.Range("I4").Value = 0.2879023314891748
$ws.Range("J4").Value = 0.2879023314891748
$ws.Range("M4").Value = 4.309709000000001
$ws.Range("N4").Value = 12.929127
$ws.Range("O4").Value = 0.170367093006711
$ws.Range("P4").Value = 0.170367093006711
$ws.Range("Q4").Value = 5.787544130341667
$ws.Range("R4").Value = 52.087897173075
$ws.Range("S4").Value = 0.04904908328566519
$ws.Range("T4").Value = 0.04904908328566519
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.342908333333333
$ws.Range("H5").Value = 4.028725
$ws.Range("I5").Value = 0.2879023314891748
$ws.Range("J5").Value = 0.2879023314891748
$ws.Range("M5").Value = 3.600126333333333
$ws.Range("N5").Value = 10.800379
$ws.Range("O5").Value = 0.1423165828288893
$ws.Range("P5").Value = 0.1423165828288893
$ws.Range("Q5").Value = 4.834639654086111
$ws.Range("R5").Value = 43.511756886775
$ws.Range("S5").Value = 0.04097327600600947
$ws.Range("T5").Value = 0.04097327600600947
$ws.Range("G6").Value = 0.4963216666666666
$ws.Range("I6").Value = 0.1064050028249084
$ws.Range("J6").Value = 0.1064050028249084
$ws.Range("M6").Value = 8.333446333333333
$ws.Range("N6").Value = 25.000339
$ws.Range("O6").Value = 0.3294294409523786
$ws.Range("P6").Value = 0.3294294409523787
$ws.Range("Q6").Value = 4.136069973237221
$ws.Range("R6").Value = 37.224629759135
$ws.Range("S6").Value = 0.03505294059514585
$ws.Range("T6").Value = 0.03505294059514585
$ws.Range("G7").Value = 0.4963216666666666
$ws.Range("I7").Value = 0.1064050028249084
$ws.Range("J7").Value = 0.1064050028249084
$ws.Range("O7").Value = 0.357886883212021
$ws.Range("P7").Value = 0.357886883212021
$ws.Range("Q7").Value = 4.493360360231666
$ws.Range("R7").Value = 40.44024324208499
$ws.Range("S7").Value = 0.03808095481917276
$ws.Range("T7").Value = 0.03808095481917276
$ws.Range("G8").Value = 0.4963216666666666
$ws.Range("I8").Value = 0.1064050028249084
$ws.Range("J8").Value = 0.1064050028249084
$ws.Range("M8").Value = 4.309709000000001
$ws.Range("N8").Value = 12.929127
$ws.Range("O8").Value = 0.170367093006711
$ws.Range("P8").Value = 0.170367093006711
$ws.Range("Q8").Value = 2.139001953728334
$ws.Range("R8").Value = 19.251017583555
$ws.Range("S8").Value = 0.01812791101265052
$ws.Range("T8").Value = 0.01812791101265052
$ws.Range("G9").Value = 0.4963216666666666
$ws.Range("I9").Value = 0.1064050028249084
$ws.Range("J9").Value = 0.1064050028249084
$ws.Range("M9").Value = 3.600126333333333
$ws.Range("N9").Value = 10.800379
$ws.Range("O9").Value = 0.1423165828288893
$ws.Range("P9").Value = 0.1423165828288893
$ws.Range("Q9").Value = 1.786820701970555
$ws.Range("R9").Value = 16.081386317735
$ws.Range("S9").Value = 0.01514319639793927
$ws.Range("T9").Value = 0.01514319639793927
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.1501973333333333
$ws.Range("H10").Value = 0.450592
$ws.Range("I10").Value = 0.03220038283833477
$ws.Range("J10").Value = 0.03220038283833477
$ws.Range("M10").Value = 8.333446333333333
$ws.Range("N10").Value = 25.000339
$ws.Range("O10").Value = 0.3294294409523786
$ws.Range("P10").Value = 0.3294294409523787
$ws.Range("Q10").Value = 1.251661416743111
$ws.Range("R10").Value = 11.264952750688
$ws.Range("S10").Value = 0.01060775411688519
$ws.Range("T10").Value = 0.01060775411688519
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.1501973333333333
$ws.Range("H11").Value = 0.450592
$ws.Range("I11").Value = 0.03220038283833477
$ws.Range("J11").Value = 0.03220038283833477
$ws.Range("O11").Value = 0.357886883212021
$ws.Range("P11").Value = 0.357886883212021
$ws.Range("Q11").Value = 1.359784972405333
$ws.Range("R11").Value = 12.238064751648
$ws.Range("S11").Value = 0.01152409465224548
$ws.Range("T11").Value = 0.01152409465224548
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.1501973333333333
$ws.Range("H12").Value = 0.450592
$ws.Range("I12").Value = 0.03220038283833477
$ws.Range("J12").Value = 0.03220038283833477
$ws.Range("M12").Value = 4.309709000000001
$ws.Range("N12").Value = 12.929127
$ws.Range("O12").Value = 0.170367093006711
$ws.Range("P12").Value = 0.170367093006711
$ws.Range("Q12").Value = 0.6473067992426668
$ws.Range("R12").Value = 5.825761193184
$ws.Range("S12").Value = 0.005485885617870282
$ws.Range("T12").Value = 0.005485885617870281
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.1501973333333333
$ws.Range("H13").Value = 0.450592
$ws.Range("I13").Value = 0.03220038283833477
$ws.Range("J13").Value = 0.03220038283833477
$ws.Range("M13").Value = 3.600126333333333
$ws.Range("N13").Value = 10.800379
$ws.Range("O13").Value = 0.1423165828288893
$ws.Range("P13").Value = 0.1423165828288893
$ws.Range("Q13").Value = 0.5407293749297777
$ws.Range("R13").Value = 4.866564374368
$ws.Range("S13").Value = 0.004582648451333814
$ws.Range("T13").Value = 0.004582648451333814
$ws.Range("G14").Value = 2.675030666666667
$ws.Range("H14").Value = 8.025092000000001
$ws.Range("I14").Value = 0.573492282847582
$ws.Range("J14").Value = 0.573492282847582
$ws.Range("M14").Value = 8.333446333333333
$ws.Range("N14").Value = 25.000339
$ws.Range("O14").Value = 0.3294294409523786
$ws.Range("P14").Value = 0.3294294409523787
$ws.Range("Q14").Value = 22.29222450068756
$ws.Range("R14").Value = 200.630020506188
$ws.Range("S14").Value = 0.1889252421289823
$ws.Range("T14").Value = 0.1889252421289824
$ws.Range("G15").Value = 2.675030666666667
$ws.Range("H15").Value = 8.025092000000001
$ws.Range("I15").Value = 0.573492282847582
$ws.Range("J15").Value = 0.573492282847582
$ws.Range("O15").Value = 0.357886883212021
$ws.Range("P15").Value = 0.357886883212021
$ws.Range("Q15").Value = 24.21791666023866
$ws.Range("R15").Value = 217.961249942148
$ws.Range("S15").Value = 0.2052453656544679
$ws.Range("T15").Value = 0.2052453656544679
$ws.Range("G16").Value = 2.675030666666667
$ws.Range("H16").Value = 8.025092000000001
$ws.Range("I16").Value = 0.573492282847582
$ws.Range("J16").Value = 0.573492282847582
$ws.Range("M16").Value = 4.309709000000001
$ws.Range("N16").Value = 12.929127
$ws.Range("O16").Value = 0.170367093006711
$ws.Range("P16").Value = 0.170367093006711
$ws.Range("Q16").Value = 11.52860373940934
$ws.Range("R16").Value = 103.757433654684
$ws.Range("S16").Value = 0.09770421309052504
$ws.Range("T16").Value = 0.09770421309052503
$ws.Range("G17").Value = 2.675030666666667
$ws.Range("H17").Value = 8.025092000000001
$ws.Range("I17").Value = 0.573492282847582
$ws.Range("J17").Value = 0.573492282847582
$ws.Range("M17").Value = 3.600126333333333
$ws.Range("N17").Value = 10.800379
$ws.Range("O17").Value = 0.1423165828288893
$ws.Range("P17").Value = 0.1423165828288893
$ws.Range("Q17").Value = 9.63044834554089
$ws.Range("R17").Value = 86.674035109868
$ws.Range("S17").Value = 0.0816174619736067
$ws.Range("T17").Value = 0.0816174619736067
